$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix recurring text typos: "tte" -> "Tte." and "mecanico" -> "mecánico" ---
$pairs = @(
    @("C8","D8", "sala muestra Tte. sub 5"),
    @("C9","D9", "oficina talleres mecánicos Tte.-sub 6"),
    @("C10","D10", "casino Tte.-sub -6"),
    @("C11","D11", "oficinas generales Tte.- sub-6 "),
    @("C17","D17", "Taller de Martillo Tte. sub-5"),
    @("C21","D21", "sala muestra Tte. sub 5"),
    @("C22","D22", "telecomando Tte. -sub-5"),
    @("C26","D26", "oficina talleres mecánicos Tte.-sub 6"),
    @("C27","D27", "casino Tte.-sub -6"),
    @("C28","D28", "oficinas generales Tte.- sub-6 "),
    @("C38","D38", "sala muestra Tte. sub 5"),
    @("C39","D39", "telecomando Tte. -sub-5"),
    @("C45","D45", "sala muestra Tte. sub 5"),
    @("C46","D46", "oficina talleres mecánicos Tte.-sub 6"),
    @("C47","D47", "casino Tte.-sub -6"),
    @("C48","D48", "oficinas generales Tte.- sub-6 "),
    @("C54","D54", "Taller de Martillo Tte. sub-5"),
    @("C58","D58", "sala muestra Tte. sub 5"),
    @("C59","D59", "telecomando Tte. -sub-5"),
    @("C65","D65", "sala muestra Tte. sub 5"),
    @("C68","D68", "oficina talleres mecánicos Tte.-sub 6"),
    @("C69","D69", "casino Tte.-sub -6"),
    @("C70","D70", "oficinas generales Tte.- sub-6 "),
    @("C75","D75", "Extractor Pozo#1 Talleres mecánico tte-sub-5"),
    @("C76","D76", "Extractor Pozo#2 Talleres mecánico tte-sub-5"),
    @("C77","D77", "Extractor Pozo#3 Talleres mecánico tte-sub-5"),
    @("C78","D78", "telecomando Tte. -sub-5"),
    @("C84","D84", "sala muestra Tte. sub 5"),
    @("C85","D85", "oficina talleres mecánicos Tte.-sub 6"),
    @("C86","D86", "casino Tte.-sub -6"),
    @("C87","D87", "oficinas generales Tte.- sub-6 ")
)

foreach ($pair in $pairs) {
    $ws.Range($pair[0]).Value = $pair[2]
    $ws.Range($pair[1]).Value = $pair[2]
}

# --- Advance the pauta one month (dates shift from Sep to Oct 2022) ---
$dates = @{
    "B5" = 44826
    "B9" = 44827
    "B14" = 44830
    "B18" = 44831
    "B22" = 44832
    "B26" = 44833
    "B31" = 44834
    "B35" = 44837
    "B39" = 44838
    "B42" = 44839
    "B46" = 44840
    "B51" = 44841
    "B55" = 44845
    "B59" = 44846
    "B62" = 44847
    "B68" = 44848
    "B73" = 44851
    "B78" = 44852
    "B81" = 44853
    "B85" = 44854
}
foreach ($ref in $dates.Keys) {
    $ws.Range($ref).Value = $dates[$ref]
}

# --- Row 42 (casino tte-7 on 44839/Oct 5): work not done, clear "Trabajos Realizados"
#     and note the reason in "Observación" ---
$ws.Range("D42").Value = ""
$ws.Range("E42").Value = "Casino en mantención eléctrica "

# --- Selection / scroll position left by the user after editing ---
$ws.Range("D26").Select()
